$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits at the very
#    start of the document, in the title paragraph). Word keeps "_GoBack"
#    hidden from Bookmarks.Count/enumeration, but it can still be reached
#    by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the contact e-mail address run and replace it with the new
#    address, split across two runs ("group-docieee2020@atha" +
#    "bascau.ca") with the "_GoBack" bookmark re-inserted at the split
#    point, exactly like Word leaves it after a human types over the
#    selection and stops there.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("cyberscitechcongress2020@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, "group-docieee2020@athabascau.ca", 2)

if ($found) {
    $splitOffset = $rng.Start + "group-docieee2020@atha".Length
    $bmRange = $d.Range($splitOffset, $splitOffset)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
